# October-2014-bank_statement.xlsx edit:
#  - remove the "awetwetawe"/BalaRaju Vankala row (3) and the Priyanka Muddana row (4)
#  - clear the account number in A2 (Sekhar Beri's row)
#  - bump Sekhar Beri's Netpay (C2) to 10000
#  - slightly narrow column A now that the long "awetwetawe" text is gone

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 3 and 4 completely (this also shrinks the used range / dimension
# from A1:D4 down to A1:D2).
$ws.Rows("3:4").Delete()

# A2 had an account number (17249172304); it becomes a blank cell.
$ws.Range("A2").ClearContents()

# Netpay for Sekhar Beri changes from 1304.86 to 10000.
$ws.Range("C2").Value = 10000.0

# Column A narrows from ~15.19 to ~14.09 (chars).
$ws.Columns("A").ColumnWidth = 13.43
